# Integreated new AUT. Added code to highlight element. Added code to take screenshot.
# Switch the TestData workbook over from the old WordPress demo site to the
# Mercury Tours demo site, and add a row used to verify the page title.

$wb = $excel.ActiveWorkbook

$wsOpenBrowser        = $wb.Worksheets.Item(1)   # openBrowser
$wsInputLoginDetails  = $wb.Worksheets.Item(2)   # inputLoginDetails
$wsVerifyEnteredVals  = $wb.Worksheets.Item(3)   # verifyEnteredValues

# ---------------------------------------------------------------------------
# Sheet "openBrowser": new URL + a new row used to check the browser title.
# ---------------------------------------------------------------------------
$wsOpenBrowser.Range("B2").Value = "http://newtours.demoaut.com/"
$wsOpenBrowser.Range("A3").Value = "title"
$wsOpenBrowser.Range("B3").Value = "Welcome: Mercury Tours"

# ---------------------------------------------------------------------------
# Sheet "inputLoginDetails": new credentials for the Mercury Tours site -
# plain values now instead of hyperlinked text.
# ---------------------------------------------------------------------------
$wsInputLoginDetails.Range("B1").Value = "mercury"
$rngLogin = $wsInputLoginDetails.Range("B2")
$rngLogin.Value = "mercury"
[void]$rngLogin.Hyperlinks.Delete()
$rngLogin.Style = "Normal"

# ---------------------------------------------------------------------------
# Sheet "verifyEnteredValues": mirrors inputLoginDetails.
# ---------------------------------------------------------------------------
$wsVerifyEnteredVals.Range("B1").Value = "mercury"
$rngVerify = $wsVerifyEnteredVals.Range("B2")
$rngVerify.Value = "mercury"
[void]$rngVerify.Hyperlinks.Delete()
$rngVerify.Style = "Normal"

# ---------------------------------------------------------------------------
# Selections / active sheet.
# Final state: openBrowser is the active (tabSelected) sheet with A3 selected;
# inputLoginDetails and verifyEnteredValues both have A3 / B2 selected
# respectively, without being the active tab.
# ---------------------------------------------------------------------------
[void]$wsInputLoginDetails.Range("A3").Select()
[void]$wsVerifyEnteredVals.Range("B2").Select()
[void]$wsOpenBrowser.Range("A3").Select()
